# Narrow both exports down to 1 month of orders (2017-11-14) instead of the
# previous 2017-12-14 range: fix the timestamp on the remaining order and
# drop the later orders' rows from both the per_product and per_order sheets.

$wb = $excel.ActiveWorkbook

$wsProduct = $wb.Worksheets.Item("per_product")
$wsOrder = $wb.Worksheets.Item("per_order")

# Correct the surviving order's timestamp on both sheets (was exported with
# the wrong month).
$wsProduct.Range("D2").Value = "2017-11-14 11:27:19"
$wsOrder.Range("D2").Value = "2017-11-14 11:27:19"

# Remove the later orders (rows 3-6) that fall outside the new date window.
$wsProduct.Range("A3:G6").Delete()
$wsOrder.Range("A3:F6").Delete()

# per_order's "items" column was sized to fit the longest combined item list;
# now that only the single-item order remains, shrink it to match (target
# OOXML width ~60.128174, closest reachable value through the COM
# character-width rounding).
$wsOrder.Columns.Item(6).ColumnWidth = 59.33
